$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3540.7
$ws.Range("I74").Value = 3547.9412
$ws.Range("K74").Value = 3547.9412
$ws.Range("M74").Value = -2611.9412
$ws.Range("H77").Value = 3540.7
$ws.Range("I77").Value = 3547.9412
$ws.Range("K77").Value = 17739.706
$ws.Range("M77").Value = -13059.706
$ws.Range("H98").Value = 971.2449
$ws.Range("I98").Value = 904.881
$ws.Range("J98").Value = 1369.4286
$ws.Range("K98").Value = 904.881
$ws.Range("L98").Value = 1369.4286
$ws.Range("M98").Value = 593.119
$ws.Range("N98").Value = -4365.4286
$ws.Range("H100").Value = 12200.6
$ws.Range("I100").Value = 7333.3335
$ws.Range("K100").Value = 7333.3335
$ws.Range("M100").Value = -6792.3335
$ws.Range("H122").Value = 971.2449
$ws.Range("I122").Value = 904.881
$ws.Range("J122").Value = 1369.4286
$ws.Range("K122").Value = 2714.643
$ws.Range("L122").Value = 4108.2858
$ws.Range("M122").Value = -264.643
$ws.Range("N122").Value = -9008.2858
$ws.Range("H132").Value = 157120.52
$ws.Range("I132").Value = 3311.0527
$ws.Range("J132").Value = 1253013
$ws.Range("K132").Value = 9933.158100000001
$ws.Range("L132").Value = 3759039
$ws.Range("M132").Value = -7403.158100000001
$ws.Range("N132").Value = -3764099
$ws.Range("H138").Value = 126615.516
$ws.Range("I138").Value = 2051.65
$ws.Range("J138").Value = 167456.12
$ws.Range("K138").Value = 6154.950000000001
$ws.Range("L138").Value = 502368.36
$ws.Range("M138").Value = -1014.950000000001
$ws.Range("N138").Value = -512648.36
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24031.969
$ws.Range("I32").Value = 15754.634
$ws.Range("J32").Value = 37101.45
$ws.Range("K32").Value = 15754.634
$ws.Range("L32").Value = 37101.45
$ws.Range("M32").Value = -15467.634
$ws.Range("N32").Value = -37675.45
$ws.Range("H61").Value = 1896.0588
$ws.Range("I61").Value = 1199.4546
$ws.Range("J61").Value = 3173.1667
$ws.Range("K61").Value = 1199.4546
$ws.Range("L61").Value = 3173.1667
$ws.Range("M61").Value = -987.4546
$ws.Range("N61").Value = -3597.1667
$ws.Range("H74").Value = 23633.4
$ws.Range("I74").Value = 26856.514
$ws.Range("J74").Value = 2683.1667
$ws.Range("K74").Value = 26856.514
$ws.Range("L74").Value = 2683.1667
$ws.Range("M74").Value = -25982.514
$ws.Range("N74").Value = -4431.1667
$ws.Range("H77").Value = 23633.4
$ws.Range("I77").Value = 26856.514
$ws.Range("J77").Value = 2683.1667
$ws.Range("K77").Value = 134282.57
$ws.Range("L77").Value = 13415.8335
$ws.Range("M77").Value = -129914.57
$ws.Range("N77").Value = -22151.8335
$ws.Range("H132").Value = 235510.2
$ws.Range("I132").Value = 35071.61
$ws.Range("J132").Value = 836826
$ws.Range("K132").Value = 105214.83
$ws.Range("L132").Value = 2510478
$ws.Range("M132").Value = -102684.83
$ws.Range("N132").Value = -2515538
$ws.Range("H136").Value = 1896.0588
$ws.Range("I136").Value = 1199.4546
$ws.Range("J136").Value = 3173.1667
$ws.Range("K136").Value = 3598.3638
$ws.Range("L136").Value = 9519.500100000001
$ws.Range("M136").Value = -1048.3638
$ws.Range("N136").Value = -14619.5001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 7915
$ws.Range("I80").Value = 20062.4
$ws.Range("J80").Value = 322.875
$ws.Range("K80").Value = 20062.4
$ws.Range("L80").Value = 322.875
$ws.Range("M80").Value = -19064.4
$ws.Range("N80").Value = -2318.875
$ws.Range("H83").Value = 7915
$ws.Range("I83").Value = 20062.4
$ws.Range("J83").Value = 322.875
$ws.Range("K83").Value = 100312
$ws.Range("L83").Value = 1614.375
$ws.Range("M83").Value = -95320
$ws.Range("N83").Value = -11598.375
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3228.8
$ws.Range("I58").Value = 960.5217
$ws.Range("J58").Value = 6297.647
$ws.Range("K58").Value = 960.5217
$ws.Range("L58").Value = 6297.647
$ws.Range("M58").Value = -757.5217
$ws.Range("N58").Value = -6703.647
$ws.Range("H132").Value = 42058.2
$ws.Range("I132").Value = 54323.367
$ws.Range("J132").Value = 3218.5
$ws.Range("K132").Value = 162970.101
$ws.Range("L132").Value = 9655.5
$ws.Range("M132").Value = -160440.101
$ws.Range("N132").Value = -14715.5
$ws.Range("H134").Value = 931.0238000000001
$ws.Range("I134").Value = 743.82355
$ws.Range("J134").Value = 1726.625
$ws.Range("K134").Value = 2231.47065
$ws.Range("L134").Value = 5179.875
$ws.Range("M134").Value = 303.5293500000002
$ws.Range("N134").Value = -10249.875
$ws.Range("H136").Value = 3228.8
$ws.Range("I136").Value = 960.5217
$ws.Range("J136").Value = 6297.647
$ws.Range("K136").Value = 2881.5651
$ws.Range("L136").Value = 18892.941
$ws.Range("M136").Value = -331.5650999999998
$ws.Range("N136").Value = -23992.941
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 132396.52
$ws.Range("J131").Value = 147919.34
$ws.Range("L131").Value = 443758.02
$ws.Range("N131").Value = -453838.02
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 68958
$ws.Range("I132").Value = 1564.75
$ws.Range("J132").Value = 93464.63
$ws.Range("K132").Value = 4694.25
$ws.Range("L132").Value = 280393.89
$ws.Range("M132").Value = -2164.25
$ws.Range("N132").Value = -285453.89
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2920.7097
$ws.Range("I7").Value = 2921.4546
$ws.Range("J7").Value = 2918.889
$ws.Range("K7").Value = 2921.4546
$ws.Range("L7").Value = 2918.889
$ws.Range("M7").Value = -2809.4546
$ws.Range("N7").Value = -3142.889
$ws.Range("H16").Value = 2990
$ws.Range("I16").Value = 2990
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2990
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2820
$ws.Range("N16").ClearContents()
$ws.Range("H46").Value = 2157.9167
$ws.Range("I46").Value = 300
$ws.Range("J46").Value = 2529.5
$ws.Range("K46").Value = 300
$ws.Range("L46").Value = 2529.5
$ws.Range("M46").Value = -112
$ws.Range("N46").Value = -2905.5
$ws.Range("H55").Value = 212.9
$ws.Range("I55").Value = 99.57143000000001
$ws.Range("J55").Value = 477.33334
$ws.Range("K55").Value = 99.57143000000001
$ws.Range("L55").Value = 477.33334
$ws.Range("M55").Value = 73.42856999999999
$ws.Range("N55").Value = -823.33334
$ws.Range("H82").Value = 1389.8422
$ws.Range("I82").Value = 1103.1666
$ws.Range("J82").Value = 1881.2858
$ws.Range("K82").Value = 1103.1666
$ws.Range("L82").Value = 1881.2858
$ws.Range("M82").Value = -742.1666
$ws.Range("N82").Value = -2603.2858
$ws.Range("H85").Value = 1389.8422
$ws.Range("I85").Value = 1103.1666
$ws.Range("J85").Value = 1881.2858
$ws.Range("K85").Value = 1103.1666
$ws.Range("L85").Value = 1881.2858
$ws.Range("M85").Value = 144.8334
$ws.Range("N85").Value = -4377.2858
$ws.Range("H126").Value = 2920.7097
$ws.Range("I126").Value = 2921.4546
$ws.Range("J126").Value = 2918.889
$ws.Range("K126").Value = 8764.363799999999
$ws.Range("L126").Value = 8756.667000000001
$ws.Range("M126").Value = -6294.363799999999
$ws.Range("N126").Value = -13696.667
$ws.Range("H132").Value = 229701.33
$ws.Range("I132").Value = 55117.58
$ws.Range("J132").Value = 671980.1
$ws.Range("K132").Value = 165352.74
$ws.Range("L132").Value = 2015940.3
$ws.Range("M132").Value = -162822.74
$ws.Range("N132").Value = -2021000.3
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1626
$ws.Range("I122").Value = 1084.8889
$ws.Range("J122").Value = 2600
$ws.Range("K122").Value = 3254.6667
$ws.Range("L122").Value = 7800
$ws.Range("M122").Value = -804.6666999999998
$ws.Range("N122").Value = -12700
$ws.Range("H126").Value = 1043.625
$ws.Range("I126").Value = 912.25
$ws.Range("J126").Value = 1175
$ws.Range("K126").Value = 2736.75
$ws.Range("L126").Value = 3525
$ws.Range("M126").Value = -266.75
$ws.Range("N126").Value = -8465
$ws.Range("H132").Value = 7042.263
$ws.Range("I132").Value = 999.4167
$ws.Range("J132").Value = 17401.428
$ws.Range("K132").Value = 2998.2501
$ws.Range("L132").Value = 52204.284
$ws.Range("M132").Value = -468.2501000000002
$ws.Range("N132").Value = -57264.284
